# Regenerate the "K" column (column G, formerly "Strike#") with freshly
# calculated strikeout values (s_vals) for each start in the save_data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values (s_vals) computed for rows 2-39, replacing the old Strike# data.
$newK = @(
    3, 7, 3, 5, 7, 5, 8, 5, 9, 5,
    8, 8, 5, 2, 6, 7, 4, 3, 3, 10,
    1, 2, 4, 4, 6, 3, 4, 6, 5, 5,
    4, 4, 7, 6, 4, 3, 4, 1
)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
